# "Generate Report for handback" -- fills in the handback report columns
# (Latest Target File / Latest Handback File / Latest Handback DateTime)
# for rows that have come back in sync with en-US, on both the zh-cn and
# de-de language sheets (and refreshes the Overview/status text).

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# HyperLink font color as used by the workbook's existing "HyperLink"
# style (RGB 6495ED), passed in BGR order expected by OLE Color values.
$hyperlinkColor = 15570276

function Style-AsHyperlink($rng) {
    $rng.Font.Underline = $true
    $rng.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------
# Overview sheet: refresh the status text shown for each source file.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("B2").Value = $newStatus
$zh.Range("B3").Value = $newStatus

# Latest Target File / Latest Handback File columns (E/F), rows 2-3,
# mirror the source file (A) and handoff file (C) for each row.
$zh.Range("E2").Value = "a.md.md"
Style-AsHyperlink $zh.Range("E2")
$zh.Range("F2").Value = "a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf"
Style-AsHyperlink $zh.Range("F2")

$zh.Range("E3").Value = "a.md.md"
Style-AsHyperlink $zh.Range("E3")
$zh.Range("F3").Value = "a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf"
Style-AsHyperlink $zh.Range("F3")

$zh.Hyperlinks.Add($zh.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/ebef2810a4f0b90d2ea0359332dbd4992cf9a247/e2e/a.md.md", "", "", "a.md.md")
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4746d5772b5ef0c17dd200de866bfd204e01c0bc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf", "", "", "a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/ebef2810a4f0b90d2ea0359332dbd4992cf9a247/e2e/a.md.md", "", "", "a.md.md")
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4746d5772b5ef0c17dd200de866bfd204e01c0bc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf", "", "", "a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf")

# Latest Handback DateTime (G), rows 2-3: the handback just completed.
$zh.Range("G2").Value = "2016-01-18 06:34:05"
$zh.Range("G3").Value = "2016-01-18 06:34:05"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("B2").Value = $newStatus
$de.Range("B3").Value = $newStatus

$de.Range("E2").Value = "a.md.md"
Style-AsHyperlink $de.Range("E2")
$de.Range("F2").Value = "a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf"
Style-AsHyperlink $de.Range("F2")

$de.Range("E3").Value = "a.md.md"
Style-AsHyperlink $de.Range("E3")
$de.Range("F3").Value = "a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf"
Style-AsHyperlink $de.Range("F3")

$de.Hyperlinks.Add($de.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/ebef2810a4f0b90d2ea0359332dbd4992cf9a247/e2e/a.md.md", "", "", "a.md.md")
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8d65e4abc890f669c226ee43b24d0e160f1f52a3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf", "", "", "a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf")
$de.Hyperlinks.Add($de.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/ebef2810a4f0b90d2ea0359332dbd4992cf9a247/e2e/a.md.md", "", "", "a.md.md")
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8d65e4abc890f669c226ee43b24d0e160f1f52a3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf", "", "", "a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf")

# Latest Handback DateTime (G), rows 2-3.
$de.Range("G2").Value = "2016-01-18 06:34:21"
$de.Range("G3").Value = "2016-01-18 06:34:21"

Write-Host "Handback report generated"
